$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 = Week 13: add the missing "content" link (column H), which was
# previously an empty placeholder cell styled differently - reset the
# style back to Normal now that it holds real content (matches H13 etc.)
$ws.Range("H14").Style = "Normal"
$ws.Range("H14").Value = "/content/12-content"

# Row 13 = Week 12: add the missing "example" link (column I)
$ws.Range("I13").Value = "/example/11-example"

# Move the active selection to I14, matching where the user ended up
$ws.Range("I14").Select()
